# Fix typos in the header row of the multimedia.txt sheet:
#   occurenceID   -> occurrenceID
#   scentificName -> scientificName
# and move the active selection to J1 (matching the saved state in Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "occurrenceID"
$ws.Range("J1").Value = "scientificName"

$ws.Range("J1").Select()
